$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: num_customers (C) 41 -> 42, retention_rate (E) recalculated as C27/D27
$ws.Range("C27").Value = 42
$ws.Range("E27").Value = 0.01865008880994671

# Row 37: num_customers (C) 731 -> 736, cohort_size (D) 731 -> 736
$ws.Range("C37").Value = 736
$ws.Range("D37").Value = 736
